$d = $word.ActiveDocument

$old = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Souhvězdí Orion 2022: 16. – 25. ledna, 14. – 23. února, 14. – 24. března"
$new = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 16. – 25. ledna, 14. – 23. února, 14. – 24. března. Při pozorování použijte hvězdy oblohy, které zobrazují souhvězdí Souhvězdí Orion.16. – 25. ledna, 14. – 23. února, 14. – 24. března"

$range = $d.Content
$range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
